$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Font.Bold = $false
    $rng.Borders.LineStyle = -4142
    $rng.Value = $value
}

Set-TextValue 'D2' '66.243.72'
Set-TextValue 'E2' '  -4.27%  '
Set-TextValue 'D3' '3.566.59'
Set-TextValue 'E3' '  -4.61%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '586.20'
Set-TextValue 'E5' '  -5.07%  '
Set-TextValue 'D6' '184.41'
Set-TextValue 'E6' '  -0.21%  '
Set-TextValue 'D7' '3.559.19'
Set-TextValue 'E7' '  -4.77%  '
Set-TextValue 'D8' '0.613'
Set-TextValue 'E8' '  -4.30%  '
Set-TextValue 'E9' '  +0.34%  '
Set-TextValue 'D10' '0.667'
Set-TextValue 'E10' '  -7.77%  '
Set-TextValue 'E11' '  -11.11%  '
Set-TextValue 'D12' '53.50'
Set-TextValue 'E12' '  -6.89%  '
Set-TextValue 'D13' '0.0000255'
Set-TextValue 'E13' '  -13.88%  '
Set-TextValue 'D14' '9.73'
Set-TextValue 'E14' '  -9.16%  '
Set-TextValue 'D15' '4.129.87'
Set-TextValue 'E15' '  -4.67%  '
Set-TextValue 'D16' '3.565.34'
Set-TextValue 'E16' '  -4.67%  '
Set-TextValue 'D17' '0.126'
Set-TextValue 'E17' '  -0.94%  '
Set-TextValue 'D18' '18.27'
Set-TextValue 'E18' '  -6.45%  '
Set-TextValue 'D19' '12.20'
Set-TextValue 'E19' '  -6.76%  '
Set-TextValue 'D20' '66.041.97'
Set-TextValue 'E20' '  -4.35%  '
Set-TextValue 'E21' '  -7.66%  '
Set-TextValue 'D22' '395.80'
Set-TextValue 'E22' '  -4.80%  '
Set-TextValue 'E23' '  -7.43%  '
Set-TextValue 'D24' '85.24'
Set-TextValue 'E24' '  -5.26%  '
Set-TextValue 'B25' 'RenderToken'
Set-TextValue 'C25' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D25' '10.91'
Set-TextValue 'E25' '  -1.00%  '
Set-TextValue 'B26' 'InternetComputer(DFINITY)'
Set-TextValue 'C26' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D26' '12.45'
Set-TextValue 'E26' '  -3.19%  '
Set-TextValue 'B27' 'ImmutableX'
Set-TextValue 'C27' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D27' '2.87'
Set-TextValue 'E27' '  -6.69%  '
Set-TextValue 'D29' '3.53'
Set-TextValue 'E29' '  -7.48%  '
Set-TextValue 'D30' '8.93'
Set-TextValue 'E30' '  -7.86%  '
Set-TextValue 'D31' '30.81'
Set-TextValue 'E31' '  -7.68%  '
Set-TextValue 'D32' '6.99'
Set-TextValue 'E32' '  -5.37%  '
Set-TextValue 'B33' 'Cosmos'
Set-TextValue 'C33' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D33' '12.12'
Set-TextValue 'E33' '  -5.48%  '
Set-TextValue 'B34' 'Bittensor'
Set-TextValue 'C34' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D34' '618.38'
Set-TextValue 'E34' '  -0.43%  '
Set-TextValue 'D35' '63.34'
Set-TextValue 'E35' '  -3.99%  '
Set-TextValue 'D36' '0.112'
Set-TextValue 'E36' '  -8.72%  '
Set-TextValue 'D37' '41.24'
Set-TextValue 'E37' '  -7.92%  '
Set-TextValue 'E38' '  -0.04%  '
Set-TextValue 'D39' '0.383'
Set-TextValue 'E39' '  -5.79%  '
Set-TextValue 'D40' '0.0₃0746'
Set-TextValue 'E40' '  -16.18%  '
Set-TextValue 'E41' '  -0.26%  '
Set-TextValue 'E42' '  -9.54%  '
Set-TextValue 'D43' '2.967.84'
Set-TextValue 'E43' '  +4.93%  '
Set-TextValue 'D44' '2.79'
Set-TextValue 'E44' '  -9.88%  '
Set-TextValue 'D45' '2.46'
Set-TextValue 'E45' '  -7.42%  '
Set-TextValue 'D46' '0.0405'
Set-TextValue 'E46' '  -9.21%  '
Set-TextValue 'B47' 'Stellar'
Set-TextValue 'C47' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D47' '0.130'
Set-TextValue 'E47' '  -5.89%  '
Set-TextValue 'B48' 'THORChain'
Set-TextValue 'C48' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D48' '8.57'
Set-TextValue 'E48' '  -7.71%  '
Set-TextValue 'B49' 'ApeXProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D49' '3.02'
Set-TextValue 'E49' '  -2.62%  '
Set-TextValue 'D50' '138.56'
Set-TextValue 'E50' '  -2.02%  '
Set-TextValue 'D51' '2.75'
Set-TextValue 'E51' '  -0.94%  '
